# Split the single "{m" / ":comment some important comment}" runs in the
# second paragraph into four runs: "{", "m", ":comment some important
# comment", "}" -- without changing the paragraph's visible text.
#
# We rebuild just that paragraph's content (runs only; the paragraph's own
# mark/properties are preserved by re-using its existing rsid attributes)
# via Range.InsertXML, which lets us place the exact run boundaries we want.

$d = $word.ActiveDocument

# Locate the paragraph that holds the "{m:comment ...}" field text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{m:comment*comment}*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing '{m:comment ... }'"
}

$range = $target.Range

$flatOpc = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979">' +
              '<w:r><w:t>{</w:t></w:r>' +
              '<w:r><w:t>m</w:t></w:r>' +
              '<w:r><w:t>:comment some important comment</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$range.InsertXML($flatOpc)
